$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-09-14 01:19:48"

# Insert a new data row at row 6 (pushes rows 6..10 down to 7..11,
# carrying their values/number formats with them; hyperlink objects stay
# anchored to their original row positions, which we fix up below).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new listing.
$ws.Cells.Item(6,1).Value = $newTimestamp
$ws.Cells.Item(6,2).Value = "【業務委託】アプリ開発の継続的パートナ募集"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5393175"
$ws.Cells.Item(6,7).Value = 93
$ws.Cells.Item(6,8).Value = "◆開発 ◇アプリ"

# Refresh the "取得日時" timestamp on every data row (2..11) to the new
# scrape time.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r,1).Value = $newTimestamp
}

# Rebuild the hyperlinks on column F (2..11) so each link's target matches
# the URL text now shown in that row.
$ws.Hyperlinks.Delete()

$urls = @(
  "https://www.lancers.jp/work/detail/5392661",
  "https://www.lancers.jp/work/detail/5392937",
  "https://www.lancers.jp/work/detail/5393052",
  "https://www.lancers.jp/work/detail/5392078",
  "https://www.lancers.jp/work/detail/5393175",
  "https://www.lancers.jp/work/detail/5392840",
  "https://www.lancers.jp/work/detail/5393015",
  "https://www.lancers.jp/work/detail/5393055",
  "https://www.lancers.jp/work/detail/5392785",
  "https://www.lancers.jp/work/detail/5392608"
)
for ($i = 0; $i -lt $urls.Length; $i++) {
    $r = $i + 2
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $urls[$i]) | Out-Null
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
